# "added std dev check"
# - Update a few Time Taken (s) values in the "General Regression" table
# - Apply a 2-decimal-place number format to the last three RMSE cells
#   (adds a new cellXfs entry, numFmtId=2 / "0.00")
# - Leave the active selection on B13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Time Taken (s)" values for the second (General Regression) table
$ws.Range("B10").Value = 8.8000000000000007
$ws.Range("B11").Value = 33.06
$ws.Range("B12").Value = 34.69

# Give the RMSE cells for rows 11-13 a 2 decimal place number format
$ws.Range("D11:D13").NumberFormat = "0.00"

# Move the active cell/selection to B13
$ws.Range("B13").Select()
